# Weekly Fruit/Vegetable price update:
# A new weekly observation is inserted as row 303 (Jengibre - Vega Modelo de
# Temuco), pushing the existing rows 303-332 down to become rows 304-333.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 303; rows 303..332 shift down to 304..333.
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with this week's data.
$ws.Cells.Item(303, 1).Value2 = 10
$ws.Cells.Item(303, 2).Value2 = 'Vega Modelo de Temuco'
$ws.Cells.Item(303, 3).Value2 = 'La Araucanía'
$ws.Cells.Item(303, 4).Value2 = 45132
$ws.Cells.Item(303, 5).Value2 = 9
$ws.Cells.Item(303, 6).Value2 = 100114007
$ws.Cells.Item(303, 7).Value2 = 'Jengibre'
$ws.Cells.Item(303, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(303, 9).Value2 = 'Primera'
$ws.Cells.Item(303, 10).Value2 = 35
$ws.Cells.Item(303, 11).Value2 = 24000
$ws.Cells.Item(303, 12).Value2 = 24000
$ws.Cells.Item(303, 13).Value2 = 24000
$ws.Cells.Item(303, 14).Value2 = '$/caja 13 kilos'
$ws.Cells.Item(303, 15).Value2 = 'Perú'
$ws.Cells.Item(303, 16).Value2 = 1846
$ws.Cells.Item(303, 17).Value2 = 13
$ws.Cells.Item(303, 18).Value2 = 'Hortaliza'
